# Commit "day 6 apache poi": cell A3 on the Employees sheet is renamed
# from "Adam" to "MadamM".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A3").Value = "MadamM"
